$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.225
    "C2" = 0.50625
    "J2" = 0.0125
    "P2" = 0.146875
    "S2" = 0.109375
    "B3" = 0.01212121212121212
    "C3" = 0.01212121212121212
    "J3" = 0.04242424242424243
    "P3" = 0.7454545454545455
    "S3" = 0.1878787878787879
    "J4" = 0.08620689655172414
    "O4" = 0.01724137931034483
    "P4" = 0.6551724137931034
    "S4" = 0.2413793103448276
    "J5" = 1
    "B6" = 0.1153846153846154
    "D6" = 0.009615384615384616
    "F6" = 0.07692307692307693
    "J6" = 0.2163461538461539
    "O6" = 0.02403846153846154
    "Q6" = 0.1682692307692308
    "R6" = 0.0625
    "S6" = 0.3269230769230769
    "B7" = 0.168141592920354
    "D7" = 0.03097345132743363
    "E7" = 0.004424778761061947
    "F7" = 0.05309734513274336
    "J7" = 0.1017699115044248
    "O7" = 0.02212389380530973
    "Q7" = 0.1504424778761062
    "R7" = 0.03982300884955752
    "S7" = 0.4292035398230089
    "B8" = 0.1096938775510204
    "D8" = 0.02551020408163265
    "E8" = 0.002551020408163265
    "F8" = 0.05612244897959184
    "J8" = 0.1020408163265306
    "O8" = 0.03571428571428571
    "Q8" = 0.1607142857142857
    "R8" = 0.07908163265306123
    "S8" = 0.4285714285714285
    "B9" = 0.09693877551020408
    "D9" = 0.02551020408163265
    "F9" = 0.07653061224489796
    "J9" = 0.1173469387755102
    "O9" = 0.01530612244897959
    "Q9" = 0.2193877551020408
    "R9" = 0.07653061224489796
    "S9" = 0.3724489795918368
    "B10" = 0.1185983827493261
    "D10" = 0.02875112309074573
    "F10" = 0.07367475292003593
    "J10" = 0.1096136567834681
    "O10" = 0.01976639712488769
    "Q10" = 0.2327044025157233
    "R10" = 0.06109613656783468
    "S10" = 0.3557951482479784
    "G11" = 0.1489971346704871
    "J11" = 0.1031518624641834
    "K11" = 0.2091690544412607
    "L11" = 0.504297994269341
    "S11" = 0.03438395415472779
    "G12" = 0.7692307692307693
    "J12" = 0.1794871794871795
    "K12" = 0.01538461538461539
    "L12" = 0.01025641025641026
    "S12" = 0.02564102564102564
    "G13" = 0.673469387755102
    "J13" = 0.2244897959183673
    "S13" = 0.1020408163265306
    "F15" = 0.02392344497607655
    "H15" = 0.1339712918660287
    "I15" = 0.1004784688995215
    "J15" = 0.3157894736842105
    "K15" = 0.05263157894736842
    "M15" = 0.02870813397129187
    "O15" = 0.07177033492822966
    "S15" = 0.2727272727272727
    "F16" = 0.0198019801980198
    "H16" = 0.1732673267326733
    "I16" = 0.06435643564356436
    "J16" = 0.4108910891089109
    "K16" = 0.1336633663366337
    "M16" = 0.02475247524752475
    "O16" = 0.06435643564356436
    "S16" = 0.1089108910891089
    "F17" = 0.0136986301369863
    "H17" = 0.1506849315068493
    "I17" = 0.1050228310502283
    "J17" = 0.3789954337899543
    "K17" = 0.1278538812785388
    "M17" = 0.0136986301369863
    "O17" = 0.0730593607305936
    "S17" = 0.136986301369863
    "F18" = 0.01428571428571429
    "H18" = 0.1928571428571429
    "I18" = 0.08571428571428572
    "J18" = 0.4357142857142857
    "K18" = 0.1
    "M18" = 0.01428571428571429
    "O18" = 0.07142857142857142
    "S18" = 0.08571428571428572
    "F19" = 0.01996672212978369
    "H19" = 0.1880199667221298
    "I19" = 0.08319467554076539
    "J19" = 0.3410981697171381
    "K19" = 0.1314475873544093
    "M19" = 0.02412645590682196
    "O19" = 0.0540765391014975
    "S19" = 0.1580698835274542
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
